$wb = $excel.ActiveWorkbook

# Sheet ALC, row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 430.08334
$ws.Range("I101").Value = 406.77777
$ws.Range("J101").Value = 500
$ws.Range("K101").Value = 1220.33331
$ws.Range("L101").Value = 1500
$ws.Range("M101").Value = 401.66669
$ws.Range("N101").Value = -4744

# Sheet ALC, row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1921706.1
$ws.Range("I127").Value = 1618.6154
$ws.Range("J127").Value = 6913933.5
$ws.Range("K127").Value = 4855.8462
$ws.Range("L127").Value = 20741800.5
$ws.Range("M127").Value = 104.1538
$ws.Range("N127").Value = -20751720.5

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5348.4287
$ws.Range("I137").Value = 5914.3335
$ws.Range("J137").Value = 4924
$ws.Range("K137").Value = 17743.0005
$ws.Range("L137").Value = 14772
$ws.Range("M137").Value = -15193.0005
$ws.Range("N137").Value = -19872

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 96976620
$ws.Range("I138").Value = 125005896
$ws.Range("J138").Value = 22231888
$ws.Range("K138").Value = 375017688
$ws.Range("L138").Value = 66695664
$ws.Range("M138").Value = -375012548
$ws.Range("N138").Value = -66705944

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2694.6667
$ws.Range("I61").Value = 1924.8235
$ws.Range("J61").Value = 4003.4
$ws.Range("K61").Value = 1924.8235
$ws.Range("L61").Value = 4003.4
$ws.Range("M61").Value = -1712.8235
$ws.Range("N61").Value = -4427.4

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2430.4
$ws.Range("I74").Value = 2015.0588
$ws.Range("J74").Value = 3313
$ws.Range("K74").Value = 2015.0588
$ws.Range("L74").Value = 3313
$ws.Range("M74").Value = -1141.0588
$ws.Range("N74").Value = -5061

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2430.4
$ws.Range("I77").Value = 2015.0588
$ws.Range("J77").Value = 3313
$ws.Range("K77").Value = 10075.294
$ws.Range("L77").Value = 16565
$ws.Range("M77").Value = -5707.294
$ws.Range("N77").Value = -25301

# Sheet ARM, row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 4694.524
$ws.Range("I97").Value = 5187.4707
$ws.Range("J97").Value = 2599.5
$ws.Range("K97").Value = 5187.4707
$ws.Range("L97").Value = 2599.5
$ws.Range("M97").Value = -4691.4707
$ws.Range("N97").Value = -3591.5

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2467.7273
$ws.Range("I132").Value = 1963.6842
$ws.Range("J132").Value = 5660
$ws.Range("K132").Value = 5891.0526
$ws.Range("L132").Value = 16980
$ws.Range("M132").Value = -3361.0526
$ws.Range("N132").Value = -22040

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2694.6667
$ws.Range("I136").Value = 1924.8235
$ws.Range("J136").Value = 4003.4
$ws.Range("K136").Value = 5774.470499999999
$ws.Range("L136").Value = 12010.2
$ws.Range("M136").Value = -3224.470499999999
$ws.Range("N136").Value = -17110.2

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2864.9575
$ws.Range("I31").Value = 2188.2
$ws.Range("J31").Value = 3366.2593
$ws.Range("K31").Value = 2188.2
$ws.Range("L31").Value = 3366.2593
$ws.Range("M31").Value = -1893.2
$ws.Range("N31").Value = -3956.2593

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2864.9575
$ws.Range("I34").Value = 2188.2
$ws.Range("J34").Value = 3366.2593
$ws.Range("K34").Value = 2188.2
$ws.Range("L34").Value = 3366.2593
$ws.Range("M34").Value = -1986.2
$ws.Range("N34").Value = -3770.2593

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6732.2593
$ws.Range("I58").Value = 7593.2856
$ws.Range("J58").Value = 5805
$ws.Range("K58").Value = 7593.2856
$ws.Range("L58").Value = 5805
$ws.Range("M58").Value = -7390.2856
$ws.Range("N58").Value = -6211

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6732.2593
$ws.Range("I136").Value = 7593.2856
$ws.Range("J136").Value = 5805
$ws.Range("K136").Value = 22779.8568
$ws.Range("L136").Value = 17415
$ws.Range("M136").Value = -20229.8568
$ws.Range("N136").Value = -22515

# Sheet CUL, row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 476.9
$ws.Range("I46").Value = 275.4
$ws.Range("J46").Value = 678.4
$ws.Range("K46").Value = 826.1999999999999
$ws.Range("L46").Value = 2035.2
$ws.Range("M46").Value = -735.1999999999999
$ws.Range("N46").Value = -2217.2

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 9760
$ws.Range("I68").Value = 2420
$ws.Range("J68").Value = 13837.777
$ws.Range("K68").Value = 7260
$ws.Range("L68").Value = 41513.331
$ws.Range("M68").Value = -6449
$ws.Range("N68").Value = -43135.331

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 9760
$ws.Range("I71").Value = 2420
$ws.Range("J71").Value = 13837.777
$ws.Range("K71").Value = 21780
$ws.Range("L71").Value = 124539.993
$ws.Range("M71").Value = -17724
$ws.Range("N71").Value = -132651.993

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1049.4
$ws.Range("I107").Value = 398.5
$ws.Range("J107").Value = 1483.3334
$ws.Range("K107").Value = 1195.5
$ws.Range("L107").Value = 4450.0002
$ws.Range("M107").Value = 724.5
$ws.Range("N107").Value = -8290.0002

# Sheet CUL, row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 6947034
$ws.Range("I121").Value = 83333780
$ws.Range("J121").Value = 2784.1365
$ws.Range("K121").Value = 250001340
$ws.Range("L121").Value = 8352.4095
$ws.Range("M121").Value = -250000030
$ws.Range("N121").Value = -10972.4095

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 675.2143
$ws.Range("I140").Value = 688.8077
$ws.Range("J140").Value = 498.5
$ws.Range("K140").Value = 2066.4231
$ws.Range("L140").Value = 1495.5
$ws.Range("M140").Value = 3113.5769
$ws.Range("N140").Value = -11855.5

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 831.7059
$ws.Range("J97").Value = 1276.8
$ws.Range("L97").Value = 1276.8
$ws.Range("N97").Value = -2268.8

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5008.3105
$ws.Range("I126").Value = 5474.2
$ws.Range("J126").Value = 4509.143
$ws.Range("K126").Value = 16422.6
$ws.Range("L126").Value = 13527.429
$ws.Range("M126").Value = -13952.6
$ws.Range("N126").Value = -18467.429

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4857.2
$ws.Range("I132").Value = 4543.6665
$ws.Range("J132").Value = 5797.8
$ws.Range("K132").Value = 13630.9995
$ws.Range("L132").Value = 17393.4
$ws.Range("M132").Value = -11100.9995
$ws.Range("N132").Value = -22453.4

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 97644.37
$ws.Range("I132").Value = 149870.28
$ws.Range("J132").Value = 6249
$ws.Range("K132").Value = 449610.84
$ws.Range("L132").Value = 18747
$ws.Range("M132").Value = -447080.84
$ws.Range("N132").Value = -23807

# Sheet WVR, row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 135249.75
$ws.Range("J14").Value = 13666.667
$ws.Range("L14").Value = 13666.667
$ws.Range("N14").Value = -14002.667

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5011
$ws.Range("J113").Value = 5375
$ws.Range("L113").Value = 16125
$ws.Range("N113").Value = -20465

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7960.7036
$ws.Range("I126").Value = 7398.4346
$ws.Range("J126").Value = 11193.75
$ws.Range("K126").Value = 22195.3038
$ws.Range("L126").Value = 33581.25
$ws.Range("M126").Value = -19725.3038
$ws.Range("N126").Value = -38521.25

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 50417.766
$ws.Range("I136").Value = 64900.273
$ws.Range("J136").Value = 23866.5
$ws.Range("K136").Value = 194700.819
$ws.Range("L136").Value = 71599.5
$ws.Range("M136").Value = -192150.819
$ws.Range("N136").Value = -76699.5
